# Fruta / hortaliza, semanal
# A new weekly price report row is inserted at row 162 (pushing the
# existing rows 162-222 down to 163-223), for the same market/product
# series (Apio, Americana (o), Primera - Macroferia Regional de Talca).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 162, shifting subsequent rows down.
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row with this week's reported values.
$ws.Cells.Item(162, 1).Value  = 5
$ws.Cells.Item(162, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(162, 3).Value  = "Maule"
$ws.Cells.Item(162, 4).Value  = 44795
$ws.Cells.Item(162, 5).Value  = 7
$ws.Cells.Item(162, 6).Value  = 100112017
$ws.Cells.Item(162, 7).Value  = "Apio"
$ws.Cells.Item(162, 8).Value  = "Americana (o)"
$ws.Cells.Item(162, 9).Value  = "Primera"
$ws.Cells.Item(162, 10).Value = 600
$ws.Cells.Item(162, 11).Value = 10000
$ws.Cells.Item(162, 12).Value = 10000
$ws.Cells.Item(162, 13).Value = 10000
$ws.Cells.Item(162, 14).Value = "$/docena de matas"
$ws.Cells.Item(162, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(162, 16).Value = 1667
$ws.Cells.Item(162, 17).Value = 6
$ws.Cells.Item(162, 18).Value = "Hortaliza"
